# Weekly update: insert a new "Acelga" price record for Vega Modelo de
# Temuco at row 576 (pushing the existing 576:605 block down to 577:606),
# matching the new dimension A1:R606.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 576:605 down to 577:606 and leave row 576 empty for the new record.
$ws.Rows.Item(576).Insert()

# Fill in the newly inserted row with the latest weekly quotation.
$ws.Range("A576").Value = 10
$ws.Range("B576").Value = "Vega Modelo de Temuco"
$ws.Range("C576").Value = "La Araucanía"
$ws.Range("D576").Value = 45267
$ws.Range("E576").Value = 9
$ws.Range("F576").Value = 100112009
$ws.Range("G576").Value = "Acelga"
$ws.Range("H576").Value = "Sin especificar"
$ws.Range("I576").Value = "Primera"
$ws.Range("J576").Value = 40
$ws.Range("K576").Value = 10000
$ws.Range("L576").Value = 10000
$ws.Range("M576").Value = 10000
$ws.Range("N576").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O576").Value = "Provincia de Cautín"
$ws.Range("P576").Value = 833
$ws.Range("Q576").Value = 12
$ws.Range("R576").Value = "Hortaliza"
